$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 0.1298123333333333
$ws.Range("H2").Value = 0.389437
$ws.Range("I2").Value = 0.01442185502613333
$ws.Range("J2").Value = 0.01442185502613333
$ws.Range("M2").Value = 123.2806423333333
$ws.Range("N2").Value = 369.841927
$ws.Range("O2").Value = 0.6241574062367528
$ws.Range("P2").Value = 0.6241574062367526
$ws.Range("Q2").Value = 16.00334783612211
$ws.Range("R2").Value = 144.030130525099
$ws.Range("S2").Value = 0.009001507626233855
$ws.Range("T2").Value = 0.009001507626233857

# Row 3
$ws.Range("G3").Value = 0.1298123333333333
$ws.Range("H3").Value = 0.389437
$ws.Range("I3").Value = 0.01442185502613333
$ws.Range("J3").Value = 0.01442185502613333
$ws.Range("O3").Value = 0.2392728888301323
$ws.Range("P3").Value = 0.2392728888301322
$ws.Range("Q3").Value = 6.134938445719444
$ws.Range("R3").Value = 55.214446011475
$ws.Range("S3").Value = 0.003450758914392284
$ws.Range("T3").Value = 0.003450758914392285

# Row 4
$ws.Range("G4").Value = 0.1298123333333333
$ws.Range("H4").Value = 0.389437
$ws.Range("I4").Value = 0.01442185502613333
$ws.Range("J4").Value = 0.01442185502613333
$ws.Range("O4").Value = 0.136569704933115
$ws.Range("P4").Value = 0.136569704933115
$ws.Range("Q4").Value = 3.501636718690445
$ws.Range("R4").Value = 31.514730468214
$ws.Range("S4").Value = 0.001969588485507191
$ws.Range("T4").Value = 0.001969588485507191

# Row 5
$ws.Range("I5").Value = 0.7859600471098795
$ws.Range("J5").Value = 0.7859600471098797
$ws.Range("M5").Value = 123.2806423333333
$ws.Range("N5").Value = 369.841927
$ws.Range("O5").Value = 0.6241574062367528
$ws.Range("P5").Value = 0.6241574062367526
$ws.Range("Q5").Value = 872.1480001291229
$ws.Range("R5").Value = 7849.332001162106
$ws.Range("S5").Value = 0.4905627844098184
$ws.Range("T5").Value = 0.4905627844098184

# Row 6
$ws.Range("I6").Value = 0.7859600471098795
$ws.Range("J6").Value = 0.7859600471098797
$ws.Range("O6").Value = 0.2392728888301323
$ws.Range("P6").Value = 0.2392728888301322
$ws.Range("S6").Value = 0.1880589309770477
$ws.Range("T6").Value = 0.1880589309770477

# Row 7
$ws.Range("I7").Value = 0.7859600471098795
$ws.Range("J7").Value = 0.7859600471098797
$ws.Range("O7").Value = 0.136569704933115
$ws.Range("P7").Value = 0.136569704933115
$ws.Range("S7").Value = 0.1073383317230134
$ws.Range("T7").Value = 0.1073383317230134

# Row 8
$ws.Range("I8").Value = 0.1996180978639869
$ws.Range("J8").Value = 0.199618097863987
$ws.Range("M8").Value = 123.2806423333333
$ws.Range("N8").Value = 369.841927
$ws.Range("O8").Value = 0.6241574062367528
$ws.Range("P8").Value = 0.6241574062367526
$ws.Range("Q8").Value = 221.5081103445918
$ws.Range("R8").Value = 1993.572993101326
$ws.Range("S8").Value = 0.1245931142007004
$ws.Range("T8").Value = 0.1245931142007004

# Row 9
$ws.Range("I9").Value = 0.1996180978639869
$ws.Range("J9").Value = 0.199618097863987
$ws.Range("O9").Value = 0.2392728888301323
$ws.Range("P9").Value = 0.2392728888301322
$ws.Range("S9").Value = 0.04776319893869221
$ws.Range("T9").Value = 0.04776319893869221

# Row 10
$ws.Range("I10").Value = 0.1996180978639869
$ws.Range("J10").Value = 0.199618097863987
$ws.Range("O10").Value = 0.136569704933115
$ws.Range("P10").Value = 0.136569704933115
$ws.Range("S10").Value = 0.02726178472459438
$ws.Range("T10").Value = 0.02726178472459437
